$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new row of data (favorite movie: Avengers)
$ws.Range("C8").Value = "avengers theme"
$ws.Range("D8").Value = "https://youtu.be/FOabQZHT4qY?t=116"
$ws.Range("E8").Value = "hope"

# Move selection to E9, matching the post-edit state
$ws.Activate()
$ws.Range("E9").Select()
